$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text format before writing so values like
# "1.005" are not auto-coerced into numbers by Excel's content-sniffing.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.475.74'
$ws.Range("E2").Value = '  -0.99%  '
$ws.Range("D3").Value = '1.831.62'
$ws.Range("E3").Value = '  -1.63%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  -0.88%  '
$ws.Range("D5").Value = '330.97'
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E7").Value = '  -2.31%  '
$ws.Range("D8").Value = '0.3835'
$ws.Range("E8").Value = '  -1.54%  '
$ws.Range("D9").Value = '46.55'
$ws.Range("E9").Value = '  -0.27%  '
$ws.Range("D10").Value = '0.07915'
$ws.Range("D11").Value = '0.9721'
$ws.Range("E11").Value = '  -3.14%  '
$ws.Range("D12").Value = '21.12'
$ws.Range("E12").Value = '  -2.18%  '
$ws.Range("D13").Value = '1.825.49'
$ws.Range("E13").Value = '  -2.59%  '
$ws.Range("D14").Value = '5.881'
$ws.Range("E14").Value = '  -1.84%  '
$ws.Range("D15").Value = '7.066'
$ws.Range("E15").Value = '  -0.88%  '
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").Value = '88.11'
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").Value = '0.06650'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").Value = '0.00001031'
$ws.Range("E19").Value = '  -1.11%  '
$ws.Range("D20").Value = '17.23'
$ws.Range("E20").Value = '  +1.89%  '
$ws.Range("D21").Value = '1.004'
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("D22").Value = '27.464.01'
$ws.Range("E22").Value = '  -1.03%  '
$ws.Range("D23").Value = '5.339'
$ws.Range("E23").Value = '  -2.19%  '
$ws.Range("D24").Value = '10.81'
$ws.Range("E24").Value = '  -0.93%  '
$ws.Range("E25").Value = '  -1.00%  '
$ws.Range("D26").Value = '2.071.52'
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("D27").Value = '157.16'
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").Value = '19.43'
$ws.Range("E28").Value = '  -1.13%  '
$ws.Range("D29").Value = '2.064'
$ws.Range("E29").Value = '  -1.06%  '
$ws.Range("D30").Value = '5.262'
$ws.Range("E30").Value = '  -2.55%  '
$ws.Range("D31").Value = '118.50'
$ws.Range("E31").Value = '  -2.02%  '
$ws.Range("D32").Value = '0.9534'
$ws.Range("E32").Value = '  -1.39%  '
$ws.Range("D33").Value = '0.09290'
$ws.Range("D34").Value = '3.570'
$ws.Range("E34").Value = '  -1.96%  '
$ws.Range("D35").Value = '5.248'
$ws.Range("E35").Value = '  -1.12%  '
$ws.Range("D36").Value = '1.319'
$ws.Range("E36").Value = '  -1.80%  '
$ws.Range("D37").Value = '0.02208'
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("D38").Value = '0.05934'
$ws.Range("D39").Value = '8.055'
$ws.Range("E39").Value = '  -0.91%  '
$ws.Range("D40").Value = '1.154'
$ws.Range("E40").Value = '  -4.36%  '
$ws.Range("D41").Value = '0.5796'
$ws.Range("E41").Value = '  -2.02%  '
$ws.Range("D42").Value = '0.1841'
$ws.Range("E42").Value = '  -2.24%  '
$ws.Range("D43").Value = '10.03'
$ws.Range("E43").Value = '  -1.77%  '
$ws.Range("D44").Value = '1.271'
$ws.Range("E44").Value = '  +0.84%  '
$ws.Range("D45").Value = '0.5495'
$ws.Range("E45").Value = '  -2.24%  '
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("D47").Value = '1.871'
$ws.Range("E47").Value = '  -2.33%  '
$ws.Range("D48").Value = '0.06642'
$ws.Range("E48").Value = '  -1.99%  '
$ws.Range("D49").Value = '110.24'
$ws.Range("E49").Value = '  -1.87%  '
$ws.Range("D50").Value = '1.040'
$ws.Range("E50").Value = '  -2.08%  '
$ws.Range("E51").Value = '  -0.91%  '

# Restore the default (unstyled) look for the Price column, matching the
# original workbook formatting.
$ws.Range("D2:D51").Style = "Normal"
